$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compras")

# Update row 3 values (customer/recipient details changed)
$ws.Range("D3").Value = "Alfred"
$ws.Range("E3").Value = 987654
$ws.Range("F3").Value = 3117889652
$ws.Range("G3").Value = "car@gmail.com"
$ws.Range("H3").Value = "Marina Sosa"
$ws.Range("I3").Value = "Itagüí"
$ws.Range("J3").Value = "Cll 45 - 78"
$ws.Range("K3").Value = "Apto"
$ws.Range("L3").Value = 3568795641
$ws.Range("M3").Value = "Cariño"
$ws.Range("N3").Value = "CAAS"
$ws.Range("Q3").Value = 15
$ws.Range("R3").Value = 3247896321

# Update row 2 Fecha value
$ws.Range("Q2").Value = 10

# Select row 4 (entire row), matching the saved selection state
$ws.Rows("4:4").Select()
